# Update the cryptocurrency price list: refresh the "Price" (D) column for
# the coins whose quote moved, and bump every "Hora" (G) cell in the data
# rows (2-51) from "5" to "6" (new hourly snapshot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing Excel to keep it as literal TEXT
# (the sheet stores numeric-looking prices/hours as text already -- a
# leading apostrophe prevents Excel from re-interpreting them as numbers),
# then strip the "number stored as text" formatting flag it adds so the
# cell style is left exactly as it was.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.ClearFormats()
}

# --- Column D: updated prices -----------------------------------------
$prices = @{
    2 = "261.54"
    3 = "22.92"
    4 = "6.201"
    5 = "0.06156"
    7 = "3.465"
    8 = "1.344"
    9 = "0.7980"
    10 = "0.1576"
    11 = "0.08101"
    12 = "0.03514"
    13 = "0.03087"
    14 = "0.09323"
    15 = "3.844"
    16 = "0.001702"
    17 = "0.04794"
    18 = "0.0006149"
    19 = "0.006200"
    20 = "0.001091"
    21 = "0.004075"
    24 = "2.205"
    25 = "0.3361"
    26 = "0.1244"
    27 = "0.0003201"
    40 = "0.04610"
    41 = "0.007099"
    42 = "0.1118"
    43 = "0.003599"
    44 = "0.009962"
    46 = "0.00005922"
    49 = "0.08947"
}
foreach ($row in $prices.Keys) {
    Set-TextValue $ws.Cells.Item($row, 4) $prices[$row]
}

# --- Column G: "Hora" snapshot marker, every data row 2-51 -------------
for ($row = 2; $row -le 51; $row++) {
    Set-TextValue $ws.Cells.Item($row, 7) "6"
}

Write-Output "Updated $($prices.Count) price cells and 50 hora cells"
